# The post for row 484 ("「アラビア語でつぶやこう」") was removed from the
# source data, so delete that entire row; Excel shifts every row below it
# up by one (old row 485 -> new row 484, ..., old row 664 -> new row 663)
# and the sheet's used-range dimension shrinks from A1:C664 to A1:C663
# automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(484).Delete()
